# diagnostic.xlsx: add a small "disconnected_elements" report to Sheet1.
#   B1 = 0            (bold, thin box border, centered/top aligned)
#   A2 = 0            (bold, thin box border, centered/top aligned)
#   B2 = "disconnected_elements"   (plain text, default formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# B1 and A2 share the same emphasised style (bold font, thin border on all
# sides, centered horizontally, top-aligned vertically). Apply the format to
# both cells together so they end up referencing a single shared cell style.
$labelCells = $excel.Union($ws.Range("B1"), $ws.Range("A2"))
foreach ($area in $labelCells.Areas) {
    $area.Font.Bold = $true
    $area.HorizontalAlignment = -4108   # xlCenter
    $area.VerticalAlignment = -4160     # xlTop
    $area.Borders.LineStyle = 1         # xlContinuous (thin box border)
}
